# Update countries & provincias Spain
# - Refresh the "last updated" timestamp banner
# - Refresh COVID stat counters for several countries (rows 4, 6, 8, 23, 27)
# - Insert "Reunion" into the sorted country list (between "Trinidad yTobago"
#   and "Letonia"), which pushes "Letonia", "Burkina Faso", "Togo" and
#   "Liberia" down by one row, and give "Reunion" its own refreshed stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Update the "last updated" banner in A1 ----------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Agosto de 2020 a las 21:56"

# ---- 2. Refresh stats for existing countries -------------------------------
# Row 4  -> Estados Unidos
$ws.Range("B4").Value = 5982555
$ws.Range("C4").Value = 26827
$ws.Range("D4").Value = 3267197
$ws.Range("E4").Value = 2532174
$ws.Range("G4").Value = 780
$ws.Range("H4").Value = 183184

# Row 6  -> Rusia
$ws.Range("B6").Value = 3307749
$ws.Range("C6").Value = 75995
$ws.Range("D6").Value = 2523443
$ws.Range("E6").Value = 723677
$ws.Range("G6").Value = 1017
$ws.Range("H6").Value = 60629

# Row 8  -> Sudafrica
$ws.Range("B8").Value = 615701
$ws.Range("C8").Value = 2684
$ws.Range("D8").Value = 525242
$ws.Range("E8").Value = 76957
$ws.Range("G8").Value = 194
$ws.Range("H8").Value = 13502

# Row 23 -> Filipinas
$ws.Range("B23").Value = 238895
$ws.Range("C23").Value = 1323
$ws.Range("E23").Value = 19944
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 9351

# Row 27 -> Bolivia
$ws.Range("B27").Value = 126225
$ws.Range("C27").Value = 256
$ws.Range("D27").Value = 112255
$ws.Range("E27").Value = 4877
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 9093

# ---- 3. Re-sort "Reunion" into its alphabetical slot -----------------------
# Currently: ... Trinidad yTobago(152), Letonia(153), Burkina Faso(154),
#            Togo(155), Liberia(156), Reunion(157), Niger(158) ...
# Target:    ... Trinidad yTobago(152), Reunion(153), Letonia(154),
#            Burkina Faso(155), Togo(156), Liberia(157), Niger(158) ...

# Remove the old "Reunion" row (157), which shifts Niger (158) etc. up by one.
$ws.Rows.Item(157).Delete()

# Insert a fresh row right after "Trinidad yTobago" (152) for "Reunion",
# pushing Letonia/Burkina Faso/Togo/Liberia back down to 154-157.
$ws.Rows.Item(153).Insert()

$ws.Range("A153").Value = "Reunion"
$ws.Range("B153").Value = 1372
$ws.Range("C153").Value = 80
$ws.Range("D153").Value = 692
$ws.Range("E153").Value = 674
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 6
